# Generate Report for Archive
# Update the localization status text from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F2/E3:F3 and the
# per-language sheets' Status column C2:C3), then narrow the now-shorter
# Status columns to match their auto-fitted content width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-language sheets: Status column (C) for both data rows
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# The Status columns were previously sized to fit "Ready for handoff"; now
# that the text is shorter, shrink them to fit "In Translation" again.
$newColumnWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth

Write-Host "Updated status text and resized Status columns on Overview, zh-cn, de-de sheets"
